{"js": "const replacements = [\n  [\"959\u00d73=2877\", \"695\u00d76=4170\"],\n  [\"687\u00d73=2061\", \"281\u00d78=2248\"],\n  [\"267\u00d79=2403\", \"226\u00d76=1356\"],\n  [\"862\u00d75=4310\", \"520\u00d79=4680\"],\n  [\"163\u00d79=1467\", \"711\u00d74=2844\"],\n  [\"829\u00d72=1658\", \"733\u00d76=4398\"],\n  [\"754\u00d79=6786\", \"853\u00d74=3412\"],\n  [\"598\u00d73=1794\", \"734\u00d78=5872\"],\n  [\"422\u00d74=1688\", \"576\u00d73=1728\"],\n  [\"432\u00d79=3888\", \"426\u00d78=3408\"],\n  [\"310\u00d74=1240\", \"721\u00d78=5768\"],\n  [\"385\u00d78=3080\", \"990\u00d79=8910\"],\n  [\"142\u00d72=284\", \"444\u00d78=3552\"],\n  [\"687\u00d74=2748\", \"371\u00d78=2968\"],\n  [\"111\u00d79=999\", \"869\u00d72=1738\"],\n  [\"196\u00d76=1176\", \"884\u00d79=7956\"],\n  [\"611\u00d73=1833\", \"842\u00d75=4210\"],\n  [\"131\u00d74=524\", \"259\u00d73=777\"],\n  [\"597\u00d77=4179\", \"315\u00d77=2205\"],\n  [\"251\u00d79=2259\", \"750\u00d77=5250\"],\n  [\"991\u00d76=5946\", \"169\u00d76=1014\"],\n  [\"601\u00d73=1803\", \"656\u00d74=2624\"],\n  [\"313\u00d76=1878\", \"489\u00d74=1956\"],\n  [\"234\u00d73=702\", \"928\u00d75=4640\"],\n  [\"868\u00d79=7812\", \"738\u00d78=5904\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"959\u00d73=2877\"; New = \"695\u00d76=4170\" },\n    @{ Old = \"687\u00d73=2061\"; New = \"281\u00d78=2248\" },\n    @{ Old = \"267\u00d79=2403\"; New = \"226\u00d76=1356\" },\n    @{ Old = \"862\u00d75=4310\"; New = \"520\u00d79=4680\" },\n    @{ Old = \"163\u00d79=1467\"; New = \"711\u00d74=2844\" },\n    @{ Old = \"829\u00d72=1658\"; New = \"733\u00d76=4398\" },\n    @{ Old = \"754\u00d79=6786\"; New = \"853\u00d74=3412\" },\n    @{ Old = \"598\u00d73=1794\"; New = \"734\u00d78=5872\" },\n    @{ Old = \"422\u00d74=1688\"; New = \"576\u00d73=1728\" },\n    @{ Old = \"432\u00d79=3888\"; New = \"426\u00d78=3408\" },\n    @{ Old = \"310\u00d74=1240\"; New = \"721\u00d78=5768\" },\n    @{ Old = \"385\u00d78=3080\"; New = \"990\u00d79=8910\" },\n    @{ Old = \"142\u00d72=284\"; New = \"444\u00d78=3552\" },\n    @{ Old = \"687\u00d74=2748\"; New = \"371\u00d78=2968\" },\n    @{ Old = \"111\u00d79=999\"; New = \"869\u00d72=1738\" },\n    @{ Old = \"196\u00d76=1176\"; New = \"884\u00d79=7956\" },\n    @{ Old = \"611\u00d73=1833\"; New = \"842\u00d75=4210\" },\n    @{ Old = \"131\u00d74=524\"; New = \"259\u00d73=777\" },\n    @{ Old = \"597\u00d77=4179\"; New = \"315\u00d77=2205\" },\n    @{ Old = \"251\u00d79=2259\"; New = \"750\u00d77=5250\" },\n    @{ Old = \"991\u00d76=5946\"; New = \"169\u00d76=1014\" },\n    @{ Old = \"601\u00d73=1803\"; New = \"656\u00d74=2624\" },\n    @{ Old = \"313\u00d76=1878\"; New = \"489\u00d74=1956\" },\n    @{ Old = \"234\u00d73=702\"; New = \"928\u00d75=4640\" },\n    @{ Old = \"868\u00d79=7812\"; New = \"738\u00d78=5904\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
